$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create "Slovakia" as a copy of "Portugal", appended at the end of the
#    workbook's sheet list.
# ---------------------------------------------------------------------------
$portugal = $wb.Worksheets.Item("Portugal")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$portugal.Copy($null, $lastSheet)
$slovakia = $wb.Worksheets.Item($wb.Worksheets.Count)
$slovakia.Name = "Slovakia"

# ---------------------------------------------------------------------------
# 2) Create "Slovakia_SlotCards215Panel" as a copy of
#    "Czech_SlotCards215Panel", appended at the end of the workbook.
# ---------------------------------------------------------------------------
$czechPanel = $wb.Worksheets.Item("Czech_SlotCards215Panel")
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$czechPanel.Copy($null, $lastSheet2)
$slovakiaPanel = $wb.Worksheets.Item($wb.Worksheets.Count)
$slovakiaPanel.Name = "Slovakia_SlotCards215Panel"

# ---------------------------------------------------------------------------
# 3) Update the market name / user-story reference text on both new sheets.
# ---------------------------------------------------------------------------
$slovakia.Range("B2").Value = "Slovakia Market"
$slovakia.Range("B4").Value = "NGC-2930/T3222"

$slovakiaPanel.Range("B2").Value = "Slovakia Market"
$slovakiaPanel.Range("B4").Value = "NGC-2930/T3222"

# ---------------------------------------------------------------------------
# 4) "Slovakia" also needs the extra "FBI800" row (like the other
#    *_SlotCards215Panel sheets) inserted right after "PCH800 5.0A" (row 8),
#    pushing "Wg" / "Slot Cards" down by one row.
# ---------------------------------------------------------------------------
$slovakia.Rows.Item(9).Insert()
$slovakia.Cells.Item(8, 1).Copy()
$slovakia.Cells.Item(9, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$slovakia.Cells.Item(9, 1).Value = "FBI800"

# ---------------------------------------------------------------------------
# 5) Selections: the previously-active sheet (Czech_SlotCards215Panel) and
#    Portugal both get their selection reset to a "select all" state, while
#    the new Slovakia / Slovakia_SlotCards215Panel sheets get the selection
#    the originals used to have.
# ---------------------------------------------------------------------------
$czechPanel.Cells.Select()

$portugal.Cells.Select()

$slovakiaPanel.Range("A8").Select()

$slovakia.Select()
$slovakia.Range("A9").Select()
